# Fix Avg_Agent_Step_Time / Avg_Experiment_Time / Std_Agent_Step_Time /
# Std_Experiment_Time columns (G, H, M, N) for rows 2-13 on Sheet1.
# Commit message: "calc correct avg exp times"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value  = 8.314492450000001
$ws.Range("H2").Value  = 448.76810782
$ws.Range("M2").Value  = 0.8883186451482461
$ws.Range("N2").Value  = 74.48029021514826

$ws.Range("G3").Value  = 9.597411920000001
$ws.Range("H3").Value  = 837.21242793
$ws.Range("M3").Value  = 1.583401387045373
$ws.Range("N3").Value  = 213.485354250579

$ws.Range("G4").Value  = 4.224720069999999
$ws.Range("H4").Value  = 134.25883026
$ws.Range("M4").Value  = 0.8253275188803025
$ws.Range("N4").Value  = 41.15293653647098

$ws.Range("G5").Value  = 4.20222397
$ws.Range("H5").Value  = 201.20463176
$ws.Range("M5").Value  = 0.6779754676739468
$ws.Range("N5").Value  = 56.9338361412727

$ws.Range("G6").Value  = 1.63645943
$ws.Range("H6").Value  = 30.86698911
$ws.Range("M6").Value  = 0.4453971724631998
$ws.Range("N6").Value  = 12.371041847014

$ws.Range("G7").Value  = 1.73958172
$ws.Range("H7").Value  = 47.02998725
$ws.Range("M7").Value  = 0.4325281065604944
$ws.Range("N7").Value  = 18.1814224254417

$ws.Range("G8").Value  = 0.83967795
$ws.Range("H8").Value  = 11.13863903
$ws.Range("M8").Value  = 0.3033605046516108
$ws.Range("N8").Value  = 5.355372929625585

$ws.Range("G9").Value  = 0.8557333499999999
$ws.Range("H9").Value  = 16.49723366
$ws.Range("M9").Value  = 0.2666376660388178
$ws.Range("N9").Value  = 7.939670751014813

$ws.Range("G10").Value = 0.43494278
$ws.Range("H10").Value = 4.3914748
$ws.Range("M10").Value = 0.1888859005646099
$ws.Range("N10").Value = 2.719136353838094

$ws.Range("G11").Value = 0.46731378
$ws.Range("H11").Value = 7.16046904
$ws.Range("M11").Value = 0.1801239939444662
$ws.Range("N11").Value = 4.705001554946197

$ws.Range("G12").Value = 0.26367343
$ws.Range("H12").Value = 2.0947919
$ws.Range("M12").Value = 0.1370732062826857
$ws.Range("N12").Value = 1.582938666465755

$ws.Range("G13").Value = 0.2694337000000001
$ws.Range("H13").Value = 3.40033725
$ws.Range("M13").Value = 0.1070338406259868
$ws.Range("N13").Value = 2.324451434122841
